$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some target values look numeric (e.g. "69.533.99", "1.00", "0.0691") or are
# percent-like text; force each target cell to Text format right before writing
# so Excel does not auto-convert the literal string into a number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.533.99'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.80%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.394.40'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +4.49%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '191.70'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '594.03'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.43%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.26%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.78'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.02%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.87%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.981.82'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +4.80%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.137'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.76'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.87%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '69.520.60'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.77%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.71%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.394.85'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +5.64%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '450.55'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +14.05%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.81'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.84'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.85'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +6.19%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("B24").Value = 'Polygon'
$ws.Range("C24").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.524'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.38%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +4.07%  '
$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.191'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.81%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.49'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.22%  '
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.19%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.01'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.93%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '23.51'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.74%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.66'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.02%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.29'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.06%  '
$ws.Range("B33").Value = 'Aptos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.00'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.24%  '
$ws.Range("B34").Value = 'USDe'
$ws.Range("C34").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.57'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +6.54%  '
$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '165.21'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.26%  '
$ws.Range("B37").Value = 'Stacks'
$ws.Range("C37").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.95'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.84%  '
$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '27.82'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +4.57%  '
$ws.Range("B39").Value = 'Mantle'
$ws.Range("C39").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.817'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +1.02%  '
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.62'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.68%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.59'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.88%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.751.89'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +5.23%  '
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.54'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.60%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '25.65'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.72%  '
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0691'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.50%  '
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.08'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.01%  '
$ws.Range("B47").Value = 'Bittensor'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '341.65'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.96%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0285'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.35%  '
$ws.Range("B49").Value = 'Arweave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '33.08'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +7.36%  '
$ws.Range("B50").Value = 'ONDO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.03'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +5.83%  '
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.36'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.42%  '
